$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2 through 66.
# Update the value from 45224 (2023-10-25) to 45233 (2023-11-03) for each row.
for ($row = 2; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value = 45233
    }
}
